# Refresh the cryptocurrency market data (price + 1h volume change) to the
# latest scrape, matching the GitHub Actions "Updated cryptos list" commit.
# Two rows (26/27 and 47/48) also swapped rank position between runs, so
# their Coin / Link / Price / Volume cells are rewritten as a full set.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: column D prices are stored as literal text in the sheet (values like
# "36.671.00" or "2.90" are not valid numbers / would lose a trailing zero).
# A leading apostrophe forces Excel to keep the entry as text instead of
# coercing it to a number, matching the original inlineStr cell content.

$ws.Range("D2").Value = "'36.671.00"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "'2.059.42"
$ws.Range("E3").Value = "  +1.31%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'243.65"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("E6").Value = "  +2.16%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'54.60"
$ws.Range("E8").Value = "  -5.12%  "
$ws.Range("D9").Value = "'58.87"
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").Value = "'0.364"
$ws.Range("E10").Value = "  -2.40%  "
$ws.Range("D11").Value = "'0.0751"
$ws.Range("E11").Value = "  -1.77%  "
$ws.Range("E12").Value = "  -2.74%  "
$ws.Range("E13").Value = "  +6.70%  "
$ws.Range("D14").Value = "'14.74"
$ws.Range("E14").Value = "  -2.90%  "
$ws.Range("D15").Value = "'2.361.88"
$ws.Range("E15").Value = "  +1.32%  "
$ws.Range("D16").Value = "'5.47"
$ws.Range("E16").Value = "  -1.76%  "
$ws.Range("D17").Value = "'2.080.72"
$ws.Range("E17").Value = "  +1.53%  "
$ws.Range("D18").Value = "'36.579.67"
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("D19").Value = "'17.04"
$ws.Range("E19").Value = "  -5.50%  "
$ws.Range("D20").Value = "'72.07"
$ws.Range("E20").Value = "  -1.49%  "
$ws.Range("D21").Value = "'0.0₃0864"
$ws.Range("E21").Value = "  -1.65%  "
$ws.Range("D22").Value = "'238.42"
$ws.Range("E22").Value = "  +1.70%  "
$ws.Range("D23").Value = "'5.26"
$ws.Range("E23").Value = "  -0.79%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("E25").Value = "  -3.18%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "'2.14"
$ws.Range("E26").Value = "  +0.66%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'9.36"
$ws.Range("E27").Value = "  -1.52%  "
$ws.Range("D28").Value = "'164.40"
$ws.Range("E28").Value = "  -1.81%  "
$ws.Range("D29").Value = "'20.16"
$ws.Range("E29").Value = "  +2.02%  "
$ws.Range("E30").Value = "  -0.83%  "
$ws.Range("D31").Value = "'5.09"
$ws.Range("E31").Value = "  -6.54%  "
$ws.Range("E32").Value = "  +8.51%  "
$ws.Range("E33").Value = "  -3.21%  "
$ws.Range("D34").Value = "'0.0597"
$ws.Range("E34").Value = "  -1.50%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("D37").Value = "'2.21"
$ws.Range("E37").Value = "  -0.31%  "
$ws.Range("D38").Value = "'0.0826"
$ws.Range("E38").Value = "  -4.26%  "
$ws.Range("E39").Value = "  -2.43%  "
$ws.Range("E40").Value = "  -4.96%  "
$ws.Range("D41").Value = "'2.90"
$ws.Range("E41").Value = "  -6.98%  "
$ws.Range("E42").Value = "  -1.82%  "
$ws.Range("E43").Value = "  -1.96%  "
$ws.Range("D44").Value = "'94.45"
$ws.Range("E44").Value = "  -1.90%  "
$ws.Range("E45").Value = "  -3.34%  "
$ws.Range("D46").Value = "'1.407.11"
$ws.Range("E46").Value = "  +9.59%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'7.61"
$ws.Range("E47").Value = "  +14.68%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'16.00"
$ws.Range("E48").Value = "  -4.41%  "
$ws.Range("D49").Value = "'2.92"
$ws.Range("E49").Value = "  +2.31%  "
$ws.Range("E50").Value = "  -2.36%  "
$ws.Range("D51").Value = "'2.250.23"
$ws.Range("E51").Value = "  +1.68%  "
